$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New log rows (20-22) ---
$ws.Range("A20").Value = 45638
$ws.Range("A20").NumberFormat = "d-mmm"
$ws.Range("B20").Value = "some animation for the background"
$ws.Range("C20").Value = 5

$ws.Range("A21").Value = 45639
$ws.Range("A21").NumberFormat = "d-mmm"
$ws.Range("B21").Value = "reworked the player aminator and added holding items as well as a background for the ice puzzle"
$ws.Range("C21").Value = 5

$ws.Range("A22").Value = 45642
$ws.Range("A22").NumberFormat = "d-mmm"
$ws.Range("B22").Value = "Worked o the aesthetics for the first ice puzzle"
$ws.Range("C22").Value = 6

# --- View/selection updates ---
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B24").Select()
